# Update "想去人数" (want-to-go count) values in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets,
# reflecting a newer scrape of the source data.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value  = 987
$wsExhibition.Range("F6").Value  = 1135
$wsExhibition.Range("F15").Value = 1368
$wsExhibition.Range("F17").Value = 1253
$wsExhibition.Range("F19").Value = 19
$wsExhibition.Range("F21").Value = 1294
$wsExhibition.Range("F26").Value = 1061
$wsExhibition.Range("F28").Value = 3267
$wsExhibition.Range("F31").Value = 1454

# --- 演出 (Performance) sheet ---
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F8").Value = 8

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 987
$wsAll.Range("F10").Value = 1135
$wsAll.Range("F18").Value = 8
$wsAll.Range("F26").Value = 1368
$wsAll.Range("F28").Value = 1253
$wsAll.Range("F30").Value = 19
$wsAll.Range("F32").Value = 1294
$wsAll.Range("F39").Value = 1061
$wsAll.Range("F41").Value = 3267
$wsAll.Range("F44").Value = 1454
